$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Text / content updates (Debtor Management + Credit Insights sections were
# reworked: new "Credit Score Descriptor Trend" breakdown replaces the old
# "Engagements" + "Credit Insights" blocks).
# ---------------------------------------------------------------------------

# Row 14 becomes the new "Debtor Management" section header (was a sub-item).
$ws.Range("A14").Value() = "Debtor Management"
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").Font.Size = 16
$ws.Range("A14").Font.ThemeColor = 2

# Row 15 becomes a regular sub-item (was the "Fee Finance" header).
$ws.Range("A15").Value() = "Number of Reminders Sent"
$ws.Range("A15").Font.Bold = $false
$ws.Range("A15").Font.Size = 14
$ws.Range("A15").Font.ThemeColor = 5

# Row 16 keeps its sub-item style, text updated.
$ws.Range("A16").Value() = "Count of Clients Dropped Credit Scores"

# Row 17 becomes a new bold sub-heading ("Credit Score Descriptor Trend").
$ws.Range("A17").Value() = "Credit Score Descriptor Trend"
$ws.Range("A17").Font.Bold = $true

# Rows 18-21 become the credit-score-descriptor breakdown, right aligned.
# (Shared-string pool order follows the order these are written: Fair, Good,
# Poor, No Data.)
$ws.Range("A19").Value() = "Fair"
$ws.Range("A19").HorizontalAlignment = -4152

$ws.Range("A18").Value() = "Good"
$ws.Range("A18").Font.Bold = $false
$ws.Range("A18").Font.Size = 14
$ws.Range("A18").Font.ThemeColor = 5
$ws.Range("A18").HorizontalAlignment = -4152

$ws.Range("A20").Value() = "Poor"
$ws.Range("A20").HorizontalAlignment = -4152

$ws.Range("A21").Value() = "No Data"
$ws.Range("A21").HorizontalAlignment = -4152

# Row 22 previously held the "Debtor Management" header - removed entirely.
$ws.Range("A22").Clear()

# Row 23 previously held "Aged Debt Trend" - text removed, style kept.
$ws.Range("A23").ClearContents()

# Row 24 previously held "Number of Reminders Sent" - removed entirely.
$ws.Range("A24").Clear()

# Row 25 previously held "Credit Insights" header - text removed, style kept.
$ws.Range("A25").ClearContents()

# Rows 26-28 previously held the old credit-score rows - removed entirely.
$ws.Range("A26").Clear()
$ws.Range("A27").Clear()
$ws.Range("A28").Clear()

# ---------------------------------------------------------------------------
# Resize the logo picture (anchor right edge shifted left by ~21.75pt).
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Width = 224.96251968503938

# ---------------------------------------------------------------------------
# Add a portrait page setup for printing / PDF export.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Update the active selection left on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()
